$wb = $excel.ActiveWorkbook

# Both mark sheets get two new columns added: "Mid Paper 1" (F) and
# "Mid Paper 2" (G), mirroring the existing "Paper 1" (D) / "Paper 2" (E)
# marks so the sheet also tracks the mid-term paper scores alongside the
# final ones (commit: "calculating percent from the total available and
# changing the A level paper grading").
$sheetNames = @("Senior Six", "Senior Five")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header row
    $ws.Cells.Item(1, 6).Value = "Mid Paper 1"
    $ws.Cells.Item(1, 7).Value = "Mid Paper 2"

    # Mirror the Paper 1 / Paper 2 marks into the new Mid Paper columns for
    # every student row.
    for ($r = 2; $r -le 6; $r++) {
        $paper1 = $ws.Cells.Item($r, 4).Value()
        $paper2 = $ws.Cells.Item($r, 5).Value()

        $ws.Cells.Item($r, 6).Value = $paper1
        $ws.Cells.Item($r, 6).NumberFormat = "General"

        $ws.Cells.Item($r, 7).Value = $paper2
        $ws.Cells.Item($r, 7).NumberFormat = "General"
    }
}

# "Senior Five" narrows its Name column and is no longer the active tab;
# selection there moves to the first new column.
$wsFive = $wb.Worksheets.Item("Senior Five")
$wsFive.Columns.Item(3).ColumnWidth = 19.6
$wsFive.Range("F1").Select()

# "Senior Six" becomes the active / selected tab, with the cursor parked
# in the newly-added data.
$wsSix = $wb.Worksheets.Item("Senior Six")
$wsSix.Activate()
$wsSix.Range("F12").Select()
